# "workers are called nodes now"
#
# On the administration slide (slide 10) there are three rectangle shapes
# labelled "Worker A", "Worker B" and "Worker C". The commit renames the
# "Worker" part of each label to "Node", leaving the trailing " A"/" B"/" C"
# untouched. Because only the first word changes, PowerPoint splits the
# paragraph into two runs: a new "Node" run (carrying the same character
# formatting as the original run) followed by the existing run whose text
# is now just " A" / " B" / " C".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

function Rename-WorkerToNode($shapeIndex) {
    $shp = $s.Shapes.Item($shapeIndex)
    $tr = $shp.TextFrame.TextRange

    # "Worker" is always the first 6 characters of the label ("Worker A/B/C").
    # Replacing just that sub-range keeps the remaining " A"/" B"/" C" text
    # (and its run formatting) intact while turning "Worker" into "Node" in
    # its own run, matching how PowerPoint itself splits edited runs.
    $worker = $tr.Characters(1, 6)
    $worker.Text = "Node"
}

# Shape indexes (within Slide 10's Shapes collection) of the three
# "Worker A" / "Worker B" / "Worker C" rectangles.
Rename-WorkerToNode 2   # Rectangle 3  -> "Worker A"
Rename-WorkerToNode 4   # Rectangle 6  -> "Worker B"
Rename-WorkerToNode 11  # Rectangle 21 -> "Worker C"
